$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(221).Insert()

for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(221, $c).Value = $ws.Cells.Item(222, $c).Value2
}

$ws.Cells.Item(221, 4).Value = 44736
$ws.Cells.Item(221, 10).Value = 500
$ws.Cells.Item(221, 11).Value = 19000
$ws.Cells.Item(221, 12).Value = 20000
$ws.Cells.Item(221, 13).Value = 19500
$ws.Cells.Item(221, 16).Value = 1950
